$d = $word.ActiveDocument

$replacements = @(
    @("84÷4=", "42÷9="),
    @("54÷4=", "30÷8="),
    @("51÷7=", "61÷6="),
    @("38÷2=", "20÷6="),
    @("99÷4=", "70÷8="),
    @("87÷8=", "17÷3="),
    @("85÷6=", "39÷3="),
    @("22÷8=", "76÷9="),
    @("50÷3=", "94÷8="),
    @("47÷9=", "92÷2="),
    @("19÷9=", "52÷2="),
    @("92÷6=", "90÷3="),
    @("13÷7=", "19÷8="),
    @("59÷3=", "17÷7="),
    @("44÷9=", "46÷2="),
    @("38÷3=", "48÷7="),
    @("10÷6=", "23÷4="),
    @("36÷5=", "92÷8="),
    @("29÷3=", "39÷9="),
    @("89÷3=", "86÷4="),
    @("74÷8=", "46÷2="),
    @("63÷7=", "43÷6="),
    @("64÷3=", "65÷4="),
    @("99÷8=", "79÷9="),
    @("77÷8=", "41÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
